# Chapter 2 tables and figures finished
# Adds a new "implantedTargets" worksheet with a 3-column table of
# implanted-target atomic composition / nitrogen-atom counts, and makes
# it the active sheet (matching the author's final workbook state).

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end of the workbook -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "implantedTargets"

# --- Header row (A1, then C1) -------------------------------------------
$ws.Range("A1").Value = 'Target'
$ws.Range("C1").Value = 'Atomic percentage'

# --- Column A: Target labels ---------------------------------------------
$ws.Range("A2").Value = 'Mo (low)'
$ws.Range("A3").Value = 'Mo (mid)'
$ws.Range("A4").Value = 'Mo (high)'
$ws.Range("A5").Value = 'Ta (low)'
$ws.Range("A6").Value = 'Ta (mid)'
$ws.Range("A7").Value = 'Ta (high)'
$ws.Range("A8").Value = 'W (low)'
$ws.Range("A9").Value = 'W (mid)'
$ws.Range("A10").Value = 'W (high)'

# --- Column C: Atomic percentage values -----------------------------------
$ws.Range("C2").Value = '11 $\pm$ 2'
$ws.Range("C3").Value = '14 $\pm$ 2'
$ws.Range("C4").Value = '26 $\pm$ 5'
$ws.Range("C5").Value = '17 $\pm$ 3'
$ws.Range("C6").Value = '26 $\pm$ 4'
$ws.Range("C7").Value = '36 $\pm$ 6'
$ws.Range("C10").Value = '22 $\pm$ 4'
$ws.Range("C9").Value = '19 $\pm$ 3'
$ws.Range("C8").Value = '13 $\pm$ 2'

# --- Column B: Nitrogen atom counts (numeric-style text, 2 decimals) ------
$ws.Range("B1").Value = 'Nitrogen atoms (10$^{17}$/cm$^{2}$)'
$ws.Range("B2:B10").NumberFormat = "0.00"
$ws.Range("B2").Value = '5.46 $\pm$ 0.11'
$ws.Range("B3").Value = '6.08 $\pm$ 0.12'
$ws.Range("B4").Value = '13.10 $\pm$ 0.66'
$ws.Range("B5").Value = '9.63 $\pm$ 0.29'
$ws.Range("B6").Value = '14.37 $\pm$ 0.57'
$ws.Range("B7").Value = '21.29 $\pm$ 1.28'
$ws.Range("B8").Value = '7.33 $\pm$ 0.15'
$ws.Range("B9").Value = '11.62 $\pm$ 0.35'
$ws.Range("B10").Value = '13.38 $\pm$ 0.54'

# --- Make the new sheet the active / selected tab ----------------------
$ws.Activate() | Out-Null
$ws.Range("K20").Select() | Out-Null
